# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values are recalculated/overwritten with updated
# strikeout (K) counts in place of the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 2
    16 = 2
    17 = 4
    18 = 1
    21 = 1
    22 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
